# Applies the "coin acceptor power / NV10 bill reader / cashless payment" edit
# to the "Mega 2560" worksheet of the PIN assignment workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mega 2560")

# --- Row 20: cashless payment (onyx) wired to an interrupt pin ---
$ws.Range("E20").Value = "cashless payment onyx"

# --- Row 19: coin acceptor power switch, with wrapped description text ---
$ws.Range("C19").Value = "coin power"
$ws.Range("E19").Value = "coin acceptor power - allows to turn off coin acceptor when all compartments are`nempty"
$ws.Range("E19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 45

# --- Row 21: NV10 bill reader wired to an interrupt pin ---
$ws.Range("E21").Value = "NV10 bill reader"

# --- New content in column E for rows 11-13 (NV10 bill reader channels / onyx inhibit) ---
$ws.Range("E11").Value = "nv 10 channel 1 open"
$ws.Range("E12").Value = "nv 10 channel 2 open"
$ws.Range("E13").Value = "onyx inhibit"

# --- Rows 20 and 21 column C: "Interrupt" pin type (reuses an existing shared string) ---
$ws.Range("C20").Value = "Interrupt"
$ws.Range("C21").Value = "Interrupt"

# --- Picture anchor: row 19 growing taller pushes the picture's computed
#     bottom-right anchor cell from row 47 up to row 45. Re-assert the
#     picture's (unchanged) visual height so the host recomputes the anchor. ---
$shape = $ws.Shapes.Item(1)
$shape.Height = 540.6824409448819

# --- Reset the view: scroll back to top and select E16 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E16").Select()
